$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("log_data")

# New data rows for 6/5, 6/6, 6/7 2020 (serial dates 43987-43989),
# continuing the daily log_total_conf series.
$dates  = @(43987, 43988, 43989)
$counts = @(3030, 3062, 3086)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 80 + $i

    # Column A: date, same display format as the existing date column (yyyy-mm-dd)
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"

    # Column B: raw daily count
    $ws.Cells.Item($row, 2).Value = $counts[$i]

    # Column C: continuation of the LOG10(B) series
    $ws.Cells.Item($row, 3).Formula = "=LOG10(B$row)"
}

# Match the author's final selection/active cell after the edit
$ws.Range("C79:C82").Select()
